$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-02-16"

# Update header label for February 2022 column (B1)
$ws.Range("B1").Value = "February 2022 (through February 16)"

# Update/add the carjacking counts that changed for 2022-02-24 data

# Englewood (row 2)
$ws.Range("B2").Value = 5

# Austin (row 3)
$ws.Range("B3").Value = 5
$ws.Range("F3").Value = 4
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 6

# South Shore (row 6)
$ws.Range("B6").Value = 4

# North Lawndale (row 8)
$ws.Range("D8").Value = 5
$ws.Range("F8").Value = 4

# Humboldt Park (row 18)
$ws.Range("B18").Value = 1
$ws.Range("L18").Value = 4

# Morgan Park (row 25)
$ws.Range("L25").Value = 3

# West Loop (row 34)
$ws.Range("J34").Value = 2

# Wicker Park (row 38)
$ws.Range("B38").Value = 2

# Rogers Park (row 39)
$ws.Range("B39").Value = 2

# River North (row 42)
$ws.Range("L42").Value = 1

# Old Town (row 43)
$ws.Range("F43").Value = 2

# Little Village (row 47)
$ws.Range("D47").Value = 1

# Woodlawn (row 57)
$ws.Range("D57").Value = 2

# Jackson Park (row 70)
$ws.Range("J70").Value = 1
